# Update scripts with new TPM (recomputed NATMI ligand-receptor edge weights).
# "Inflammatory-Mac" cluster is renamed to "Resolving-Mac" and the full
# 4x4 sending x target cluster grid (ECs, FAPs, MuSCs, Resolving-Mac) is
# now populated, including the four rows that involve Resolving-Mac as a
# sending cluster (rows 14-17), which were previously absent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vegfc"
$ws.Cells.Item(2, 3).Value = "Vipr2"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3.0
$ws.Cells.Item(2, 6).Value = 1.0
$ws.Cells.Item(2, 7).Value = 4.948739333333333
$ws.Cells.Item(2, 8).Value = 14.846218
$ws.Cells.Item(2, 9).Value = 0.4917593264632457
$ws.Cells.Item(2, 10).Value = 0.4917593264632457
$ws.Cells.Item(2, 11).Value = 1.0
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.1382803333333333
$ws.Cells.Item(2, 14).Value = 0.414841
$ws.Cells.Item(2, 15).Value = 0.03059328965493693
$ws.Cells.Item(2, 16).Value = 0.03059328965493693
$ws.Cells.Item(2, 17).Value = 0.6843133245931111
$ws.Cells.Item(2, 18).Value = 6.158819921338
$ws.Cells.Item(2, 19).Value = 0.01504453551500677
$ws.Cells.Item(2, 20).Value = 0.01504453551500677

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vegfc"
$ws.Cells.Item(3, 3).Value = "Vipr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3.0
$ws.Cells.Item(3, 6).Value = 1.0
$ws.Cells.Item(3, 7).Value = 4.948739333333333
$ws.Cells.Item(3, 8).Value = 14.846218
$ws.Cells.Item(3, 9).Value = 0.4917593264632457
$ws.Cells.Item(3, 10).Value = 0.4917593264632457
$ws.Cells.Item(3, 11).Value = 3.0
$ws.Cells.Item(3, 12).Value = 1.0
$ws.Cells.Item(3, 13).Value = 0.696771
$ws.Cells.Item(3, 14).Value = 2.090313
$ws.Cells.Item(3, 15).Value = 0.1541543653555945
$ws.Cells.Item(3, 16).Value = 0.1541543653555945
$ws.Cells.Item(3, 17).Value = 3.448138054026
$ws.Cells.Item(3, 18).Value = 31.033242486234
$ws.Cells.Item(3, 19).Value = 0.07580684687863626
$ws.Cells.Item(3, 20).Value = 0.07580684687863627

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vegfc"
$ws.Cells.Item(4, 3).Value = "Vipr2"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3.0
$ws.Cells.Item(4, 6).Value = 1.0
$ws.Cells.Item(4, 7).Value = 4.948739333333333
$ws.Cells.Item(4, 8).Value = 14.846218
$ws.Cells.Item(4, 9).Value = 0.4917593264632457
$ws.Cells.Item(4, 10).Value = 0.4917593264632457
$ws.Cells.Item(4, 11).Value = 3.0
$ws.Cells.Item(4, 12).Value = 1.0
$ws.Cells.Item(4, 13).Value = 3.682798
$ws.Cells.Item(4, 14).Value = 11.048394
$ws.Cells.Item(4, 15).Value = 0.8147861900435764
$ws.Cells.Item(4, 16).Value = 0.8147861900435764
$ws.Cells.Item(4, 17).Value = 18.22520731932133
$ws.Cells.Item(4, 18).Value = 164.026865873892
$ws.Cells.Item(4, 19).Value = 0.4006787080273832
$ws.Cells.Item(4, 20).Value = 0.4006787080273833

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Vegfc"
$ws.Cells.Item(5, 3).Value = "Vipr2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3.0
$ws.Cells.Item(5, 6).Value = 1.0
$ws.Cells.Item(5, 7).Value = 4.948739333333333
$ws.Cells.Item(5, 8).Value = 14.846218
$ws.Cells.Item(5, 9).Value = 0.4917593264632457
$ws.Cells.Item(5, 10).Value = 0.4917593264632457
$ws.Cells.Item(5, 11).Value = 1.0
$ws.Cells.Item(5, 12).Value = 0.3333333333333333
$ws.Cells.Item(5, 13).Value = 0.002107
$ws.Cells.Item(5, 14).Value = 0.006321
$ws.Cells.Item(5, 15).Value = 0.0004661549458921764
$ws.Cells.Item(5, 16).Value = 0.0004661549458921764
$ws.Cells.Item(5, 17).Value = 0.01042699377533333
$ws.Cells.Item(5, 18).Value = 0.09384294397800001
$ws.Cells.Item(5, 19).Value = 0.0002292360422194474
$ws.Cells.Item(5, 20).Value = 0.0002292360422194474

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vegfc"
$ws.Cells.Item(6, 3).Value = "Vipr2"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3.0
$ws.Cells.Item(6, 6).Value = 1.0
$ws.Cells.Item(6, 7).Value = 3.979395333333333
$ws.Cells.Item(6, 8).Value = 11.938186
$ws.Cells.Item(6, 9).Value = 0.395435006178203
$ws.Cells.Item(6, 10).Value = 0.395435006178203
$ws.Cells.Item(6, 11).Value = 1.0
$ws.Cells.Item(6, 12).Value = 0.3333333333333333
$ws.Cells.Item(6, 13).Value = 0.1382803333333333
$ws.Cells.Item(6, 14).Value = 0.414841
$ws.Cells.Item(6, 15).Value = 0.03059328965493693
$ws.Cells.Item(6, 16).Value = 0.03059328965493693
$ws.Cells.Item(6, 17).Value = 0.5502721131584445
$ws.Cells.Item(6, 18).Value = 4.952449018426
$ws.Cells.Item(6, 19).Value = 0.01209765768371154
$ws.Cells.Item(6, 20).Value = 0.01209765768371154

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vegfc"
$ws.Cells.Item(7, 3).Value = "Vipr2"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3.0
$ws.Cells.Item(7, 6).Value = 1.0
$ws.Cells.Item(7, 7).Value = 3.979395333333333
$ws.Cells.Item(7, 8).Value = 11.938186
$ws.Cells.Item(7, 9).Value = 0.395435006178203
$ws.Cells.Item(7, 10).Value = 0.395435006178203
$ws.Cells.Item(7, 11).Value = 3.0
$ws.Cells.Item(7, 12).Value = 1.0
$ws.Cells.Item(7, 13).Value = 0.696771
$ws.Cells.Item(7, 14).Value = 2.090313
$ws.Cells.Item(7, 15).Value = 0.1541543653555945
$ws.Cells.Item(7, 16).Value = 0.1541543653555945
$ws.Cells.Item(7, 17).Value = 2.772727265802
$ws.Cells.Item(7, 18).Value = 24.954545392218
$ws.Cells.Item(7, 19).Value = 0.06095803241678649
$ws.Cells.Item(7, 20).Value = 0.0609580324167865

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Vegfc"
$ws.Cells.Item(8, 3).Value = "Vipr2"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3.0
$ws.Cells.Item(8, 6).Value = 1.0
$ws.Cells.Item(8, 7).Value = 3.979395333333333
$ws.Cells.Item(8, 8).Value = 11.938186
$ws.Cells.Item(8, 9).Value = 0.395435006178203
$ws.Cells.Item(8, 10).Value = 0.395435006178203
$ws.Cells.Item(8, 11).Value = 3.0
$ws.Cells.Item(8, 12).Value = 1.0
$ws.Cells.Item(8, 13).Value = 3.682798
$ws.Cells.Item(8, 14).Value = 11.048394
$ws.Cells.Item(8, 15).Value = 0.8147861900435764
$ws.Cells.Item(8, 16).Value = 0.8147861900435764
$ws.Cells.Item(8, 17).Value = 14.65530917480933
$ws.Cells.Item(8, 18).Value = 131.897782573284
$ws.Cells.Item(8, 19).Value = 0.3221949820937961
$ws.Cells.Item(8, 20).Value = 0.3221949820937962

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Vegfc"
$ws.Cells.Item(9, 3).Value = "Vipr2"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3.0
$ws.Cells.Item(9, 6).Value = 1.0
$ws.Cells.Item(9, 7).Value = 3.979395333333333
$ws.Cells.Item(9, 8).Value = 11.938186
$ws.Cells.Item(9, 9).Value = 0.395435006178203
$ws.Cells.Item(9, 10).Value = 0.395435006178203
$ws.Cells.Item(9, 11).Value = 1.0
$ws.Cells.Item(9, 12).Value = 0.3333333333333333
$ws.Cells.Item(9, 13).Value = 0.002107
$ws.Cells.Item(9, 14).Value = 0.006321
$ws.Cells.Item(9, 15).Value = 0.0004661549458921764
$ws.Cells.Item(9, 16).Value = 0.0004661549458921764
$ws.Cells.Item(9, 17).Value = 0.008384585967333333
$ws.Cells.Item(9, 18).Value = 0.075461273706
$ws.Cells.Item(9, 19).Value = 0.0001843339839088727
$ws.Cells.Item(9, 20).Value = 0.0001843339839088727

$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Vegfc"
$ws.Cells.Item(10, 3).Value = "Vipr2"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3.0
$ws.Cells.Item(10, 6).Value = 1.0
$ws.Cells.Item(10, 7).Value = 1.113241666666667
$ws.Cells.Item(10, 8).Value = 3.339725
$ws.Cells.Item(10, 9).Value = 0.1106235215306998
$ws.Cells.Item(10, 10).Value = 0.1106235215306998
$ws.Cells.Item(10, 11).Value = 1.0
$ws.Cells.Item(10, 12).Value = 0.3333333333333333
$ws.Cells.Item(10, 13).Value = 0.1382803333333333
$ws.Cells.Item(10, 14).Value = 0.414841
$ws.Cells.Item(10, 15).Value = 0.03059328965493693
$ws.Cells.Item(10, 16).Value = 0.03059328965493693
$ws.Cells.Item(10, 17).Value = 0.1539394287472223
$ws.Cells.Item(10, 18).Value = 1.385454858725
$ws.Cells.Item(10, 19).Value = 0.003384337436837852
$ws.Cells.Item(10, 20).Value = 0.003384337436837852

$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Vegfc"
$ws.Cells.Item(11, 3).Value = "Vipr2"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3.0
$ws.Cells.Item(11, 6).Value = 1.0
$ws.Cells.Item(11, 7).Value = 1.113241666666667
$ws.Cells.Item(11, 8).Value = 3.339725
$ws.Cells.Item(11, 9).Value = 0.1106235215306998
$ws.Cells.Item(11, 10).Value = 0.1106235215306998
$ws.Cells.Item(11, 11).Value = 3.0
$ws.Cells.Item(11, 12).Value = 1.0
$ws.Cells.Item(11, 13).Value = 0.696771
$ws.Cells.Item(11, 14).Value = 2.090313
$ws.Cells.Item(11, 15).Value = 0.1541543653555945
$ws.Cells.Item(11, 16).Value = 0.1541543653555945
$ws.Cells.Item(11, 17).Value = 0.7756745093250001
$ws.Cells.Item(11, 18).Value = 6.981070583925002
$ws.Cells.Item(11, 19).Value = 0.01705309875496598
$ws.Cells.Item(11, 20).Value = 0.01705309875496598

$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Vegfc"
$ws.Cells.Item(12, 3).Value = "Vipr2"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3.0
$ws.Cells.Item(12, 6).Value = 1.0
$ws.Cells.Item(12, 7).Value = 1.113241666666667
$ws.Cells.Item(12, 8).Value = 3.339725
$ws.Cells.Item(12, 9).Value = 0.1106235215306998
$ws.Cells.Item(12, 10).Value = 0.1106235215306998
$ws.Cells.Item(12, 11).Value = 3.0
$ws.Cells.Item(12, 12).Value = 1.0
$ws.Cells.Item(12, 13).Value = 3.682798
$ws.Cells.Item(12, 14).Value = 11.048394
$ws.Cells.Item(12, 15).Value = 0.8147861900435764
$ws.Cells.Item(12, 16).Value = 0.8147861900435764
$ws.Cells.Item(12, 17).Value = 4.099844183516668
$ws.Cells.Item(12, 18).Value = 36.89859765165001
$ws.Cells.Item(12, 19).Value = 0.09013451763720247
$ws.Cells.Item(12, 20).Value = 0.09013451763720247

$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Vegfc"
$ws.Cells.Item(13, 3).Value = "Vipr2"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3.0
$ws.Cells.Item(13, 6).Value = 1.0
$ws.Cells.Item(13, 7).Value = 1.113241666666667
$ws.Cells.Item(13, 8).Value = 3.339725
$ws.Cells.Item(13, 9).Value = 0.1106235215306998
$ws.Cells.Item(13, 10).Value = 0.1106235215306998
$ws.Cells.Item(13, 11).Value = 1.0
$ws.Cells.Item(13, 12).Value = 0.3333333333333333
$ws.Cells.Item(13, 13).Value = 0.002107
$ws.Cells.Item(13, 14).Value = 0.006321
$ws.Cells.Item(13, 15).Value = 0.0004661549458921764
$ws.Cells.Item(13, 16).Value = 0.0004661549458921764
$ws.Cells.Item(13, 17).Value = 0.002345600191666667
$ws.Cells.Item(13, 18).Value = 0.02111040172500001
$ws.Cells.Item(13, 19).Value = 0.00005156770169354539
$ws.Cells.Item(13, 20).Value = 0.00005156770169354539

$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Vegfc"
$ws.Cells.Item(14, 3).Value = "Vipr2"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1.0
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.02195966666666667
$ws.Cells.Item(14, 8).Value = 0.065879
$ws.Cells.Item(14, 9).Value = 0.002182145827851387
$ws.Cells.Item(14, 10).Value = 0.002182145827851387
$ws.Cells.Item(14, 11).Value = 1.0
$ws.Cells.Item(14, 12).Value = 0.3333333333333333
$ws.Cells.Item(14, 13).Value = 0.1382803333333333
$ws.Cells.Item(14, 14).Value = 0.414841
$ws.Cells.Item(14, 15).Value = 0.03059328965493693
$ws.Cells.Item(14, 16).Value = 0.03059328965493693
$ws.Cells.Item(14, 17).Value = 0.003036590026555556
$ws.Cells.Item(14, 18).Value = 0.027329310239
$ws.Cells.Item(14, 19).Value = 0.00006675901938076962
$ws.Cells.Item(14, 20).Value = 0.00006675901938076962

$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Vegfc"
$ws.Cells.Item(15, 3).Value = "Vipr2"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1.0
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.02195966666666667
$ws.Cells.Item(15, 8).Value = 0.065879
$ws.Cells.Item(15, 9).Value = 0.002182145827851387
$ws.Cells.Item(15, 10).Value = 0.002182145827851387
$ws.Cells.Item(15, 11).Value = 3.0
$ws.Cells.Item(15, 12).Value = 1.0
$ws.Cells.Item(15, 13).Value = 0.696771
$ws.Cells.Item(15, 14).Value = 2.090313
$ws.Cells.Item(15, 15).Value = 0.1541543653555945
$ws.Cells.Item(15, 16).Value = 0.1541543653555945
$ws.Cells.Item(15, 17).Value = 0.015300858903
$ws.Cells.Item(15, 18).Value = 0.137707730127
$ws.Cells.Item(15, 19).Value = 0.000336387305205789
$ws.Cells.Item(15, 20).Value = 0.000336387305205789

$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Vegfc"
$ws.Cells.Item(16, 3).Value = "Vipr2"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1.0
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.02195966666666667
$ws.Cells.Item(16, 8).Value = 0.065879
$ws.Cells.Item(16, 9).Value = 0.002182145827851387
$ws.Cells.Item(16, 10).Value = 0.002182145827851387
$ws.Cells.Item(16, 11).Value = 3.0
$ws.Cells.Item(16, 12).Value = 1.0
$ws.Cells.Item(16, 13).Value = 3.682798
$ws.Cells.Item(16, 14).Value = 11.048394
$ws.Cells.Item(16, 15).Value = 0.8147861900435764
$ws.Cells.Item(16, 16).Value = 0.8147861900435764
$ws.Cells.Item(16, 17).Value = 0.08087301648066667
$ws.Cells.Item(16, 18).Value = 0.7278571483259999
$ws.Cells.Item(16, 19).Value = 0.001777982285194518
$ws.Cells.Item(16, 20).Value = 0.001777982285194518

$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Vegfc"
$ws.Cells.Item(17, 3).Value = "Vipr2"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1.0
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.02195966666666667
$ws.Cells.Item(17, 8).Value = 0.065879
$ws.Cells.Item(17, 9).Value = 0.002182145827851387
$ws.Cells.Item(17, 10).Value = 0.002182145827851387
$ws.Cells.Item(17, 11).Value = 1.0
$ws.Cells.Item(17, 12).Value = 0.3333333333333333
$ws.Cells.Item(17, 13).Value = 0.002107
$ws.Cells.Item(17, 14).Value = 0.006321
$ws.Cells.Item(17, 15).Value = 0.0004661549458921764
$ws.Cells.Item(17, 16).Value = 0.0004661549458921764
$ws.Cells.Item(17, 17).Value = 0.00004626901766666666
$ws.Cells.Item(17, 18).Value = 0.000416421159
$ws.Cells.Item(17, 19).Value = 0.000001017218070310902
$ws.Cells.Item(17, 20).Value = 0.000001017218070310902
